# Applies the "Updated symbol list" price/volume refresh for Thu Jan 19 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D (Price) / E (Volume(1h)) hold plain text in the source feed (e.g. "293.96",
# "1.42%") rather than numeric/percentage values, so force each touched cell to
# Text format *before* writing it -- otherwise Excel auto-converts the literal
# into a real number (dropping the trailing zeros) or a percentage (0.0142).
# NumberFormat must be set per-cell: assigning it on a multi-area (comma-joined)
# Range only takes effect on the first area.

$updates = @{
    "D2" = "293.96"
    "E2" = "1.42%"
    "E3" = "0.52%"
    "D4" = "4.937"
    "E4" = "1.42%"
    "D5" = "0.07349"
    "E5" = "2.72%"
    "D6" = "2.255"
    "E6" = "23.08%"
    "E7" = "0.69%"
    "D8" = "3.772"
    "E8" = "0.11%"
    "D9" = "0.9120"
    "E9" = "2.06%"
    "D10" = "0.1683"
    "E10" = "2.56%"
    "D11" = "0.08145"
    "E11" = "8.12%"
    "D12" = "0.08259"
    "E12" = "1.93%"
    "D13" = "0.03118"
    "E13" = "4.48%"
    "D14" = "0.1008"
    "E14" = "0.87%"
    "D15" = "0.001512"
    "E15" = "1.44%"
    "D16" = "0.005816"
    "E16" = "-0.46%"
    "D17" = "3.483"
    "E17" = "0.42%"
    "E18" = "-1.19%"
    "D19" = "0.3328"
    "E19" = "1.53%"
    "D20" = "0.1304"
    "E20" = "0.91%"
    "D21" = "3.969"
    "E21" = "-6.94%"
    "E22" = "5.09%"
    "E23" = "1.67%"
    "D24" = "0.001210"
    "E24" = "-0.16%"
    "D25" = "0.004337"
    "E25" = "-7.00%"
    "D26" = "0.0001300"
    "E26" = "4.00%"
    "D27" = "0.0003396"
    "D39" = "0.01601"
    "E39" = "-2.23%"
    "E40" = "1.72%"
    "D41" = "0.007344"
    "E41" = "-0.53%"
    "D42" = "0.008783"
    "D43" = "0.1325"
    "E43" = "1.26%"
    "D44" = "0.002044"
    "E44" = "2.06%"
    "D45" = "0.009124"
    "E45" = "-11.03%"
    "D46" = "0.00005962"
    "E46" = "2.02%"
    "E47" = "0.02%"
    "E48" = "1.41%"
    "D50" = "0.00002101"
    "E50" = "0.02%"
    "D51" = "0.0002001"
    "E51" = "0.02%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
